# This script reproduces a stock-count correction pass over the
# "Companywise Stock Value Details" report. For each affected item row,
# only the Qty (F) and its dependent Value (G, = Rate * Qty) are corrected;
# each company's "Sub Total:" row (column B) is then corrected to match the
# new sum of its rows' Value column, and the two report-wide Grand Total
# rows are corrected to match the new sum of all sub-totals.
#
# A few rows (126/127, 161/162, 290/291) are full row swaps: two adjacent
# line items traded places (their Code/Name/Rate/MRP/Qty/Value), while the
# serial number (column A) stayed put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ASQUARE FOOD BEVERAGES PRIVATE LIMITED ---
$ws.Range("F36").Value = 106
$ws.Range("G36").Value = 20857.62
$ws.Range("F45").Value = 88
$ws.Range("G45").Value = 2001.12
$ws.Range("F51").Value = 158
$ws.Range("G51").Value = 14779.32
$ws.Range("F55").Value = 135
$ws.Range("G55").Value = 7527.6
$ws.Range("F61").Value = 249
$ws.Range("G61").Value = 64921.77
$ws.Range("B66").Value = 225802.41
# --- BAJAJ ELECTRICALS LIMITED ---
$ws.Range("F78").Value = 113
$ws.Range("G78").Value = 3539.16
$ws.Range("B83").Value = 16172.28
# --- BHAWAR SALES CORPORATION ---
$ws.Range("F103").Value = 22
$ws.Range("G103").Value = 1560.24
$ws.Range("B123").Value = 76681.78
# --- BLUE STAR LIMITED ---
$ws.Range("B126").Value = 64196
$ws.Range("B127").Value = 65258
# --- COLGATE PALMOLIVE INDIA LTD ---
$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88
# --- DABUR INDIA LIMITED ---
$ws.Range("F184").Value = 65
$ws.Range("G184").Value = 5330
$ws.Range("B193").Value = 70801.17999999999
# --- Glaxosmithkline Asia Private Limited ---
$ws.Range("F212").Value = 72
$ws.Range("G212").Value = 6414.48
$ws.Range("F215").Value = 181
$ws.Range("G215").Value = 20324.49
$ws.Range("B218").Value = 85028.3
# --- GODREJ CONSUMER PRODUCTS LIMITED ---
$ws.Range("F222").Value = 1140
$ws.Range("G222").Value = 21090
$ws.Range("F223").Value = 43
$ws.Range("G223").Value = 921.49
$ws.Range("F227").Value = 57
$ws.Range("G227").Value = 6532.2
$ws.Range("B229").Value = 35188.57
# --- HIMALAYA WELLNESS COMPANY ---
$ws.Range("F267").Value = 140
$ws.Range("G267").Value = 5947.2
$ws.Range("F282").Value = 19
$ws.Range("G282").Value = 1833.12
$ws.Range("F287").Value = 71
$ws.Range("G287").Value = 3886.54
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 27
$ws.Range("G291").Value = 2313.36
$ws.Range("B295").Value = 133624.98
# --- HINDUSTAN UNILIVER LTD ---
$ws.Range("F309").Value = 80
$ws.Range("G309").Value = 3811.2
$ws.Range("F325").Value = 50
$ws.Range("G325").Value = 7558.5
$ws.Range("B328").Value = 4682.65
# --- Kanav Enterprises Private Limited ---
$ws.Range("F358").Value = 47
$ws.Range("G358").Value = 10820.34
$ws.Range("F361").Value = 264
$ws.Range("G361").Value = 37115.76
$ws.Range("B363").Value = 82369.88
# --- KARNATAKA SOAPS & DETERGENTS LTD ---
$ws.Range("F370").Value = 240
$ws.Range("G370").Value = 39837.6
$ws.Range("B372").Value = 67365.91
# --- LIFE STYLE FOODS PVT LTD ---
$ws.Range("F396").Value = 141
$ws.Range("G396").Value = 3592.68
$ws.Range("B417").Value = 182062.91
# --- N.RANGA RAO & SONS PVT LTD ---
$ws.Range("F430").Value = 241
$ws.Range("G430").Value = 11153.48
$ws.Range("F433").Value = 155
$ws.Range("G433").Value = 1494.2
$ws.Range("F437").Value = 4
$ws.Range("G437").Value = 193.44
$ws.Range("B438").Value = 27586.54
# --- PNB Kitchenware (Sridhi Enterprises Pvt Ltd) ---
$ws.Range("F471").Value = 17
$ws.Range("G471").Value = 8224.77
$ws.Range("B476").Value = 54145.21
# --- RECKITT BENCKISER INDIA PVT LTD ---
$ws.Range("F519").Value = 427
$ws.Range("G519").Value = 23433.76
$ws.Range("B525").Value = 133007
# --- SARATHI INTERNATIONAL INC ---
$ws.Range("F528").Value = 309
$ws.Range("G528").Value = 4900.74
$ws.Range("F529").Value = 130
$ws.Range("G529").Value = 4304.3
$ws.Range("F530").Value = 32
$ws.Range("G530").Value = 1381.76
$ws.Range("F531").Value = 228
$ws.Range("G531").Value = 7549.08
$ws.Range("F532").Value = 15
$ws.Range("G532").Value = 647.7
$ws.Range("F533").Value = 7
$ws.Range("G533").Value = 210.56
$ws.Range("F534").Value = 141
$ws.Range("G534").Value = 6170.16
$ws.Range("B535").Value = 27217.12
# --- Shree Raghavendra Enterprises (TTK Prestige) ---
$ws.Range("F543").Value = 5
$ws.Range("G543").Value = 1851.35
$ws.Range("B556").Value = 64487.16
# --- SOUTHERN HEALTH FOODS PVT LTD ---
$ws.Range("F558").Value = 227
$ws.Range("G558").Value = 27659.95
$ws.Range("B561").Value = 32712.81
# --- Tip Top Food Tech India ---
$ws.Range("F609").Value = 35
$ws.Range("G609").Value = 3808.35
$ws.Range("F615").Value = 108
$ws.Range("G615").Value = 16704.36
$ws.Range("F620").Value = 380
$ws.Range("G620").Value = 29864.2
$ws.Range("F622").Value = 499
$ws.Range("G622").Value = 51352.09
$ws.Range("F625").Value = 346
$ws.Range("G625").Value = 12743.18
$ws.Range("B628").Value = 227019.51
# --- VVD AND SONS PRIVATE LIMITED ---
$ws.Range("F673").Value = 0
$ws.Range("G673").Value = 0
$ws.Range("F674").Value = 975
$ws.Range("G674").Value = 159032.25
$ws.Range("F675").Value = 0
$ws.Range("G675").Value = 0
$ws.Range("F676").Value = 0
$ws.Range("G676").Value = 0
$ws.Range("F677").Value = 7
$ws.Range("G677").Value = 1012.55
$ws.Range("F678").Value = 0
$ws.Range("G678").Value = 0
$ws.Range("F679").Value = 0
$ws.Range("G679").Value = 0
$ws.Range("B680").Value = 160044.8
# --- WIPRO ENTERPRISES PVT LTD ---
$ws.Range("F687").Value = 1
$ws.Range("G687").Value = 40.99
$ws.Range("B691").Value = 11395.88
# --- XO FOOTWEAR PVT LTD ---
$ws.Range("F711").Value = 16
$ws.Range("G711").Value = 8541.6
$ws.Range("B713").Value = 73001.74000000001
# --- ZYDUS WELLNESS PRODUCTS LTD (grand totals) ---
$ws.Range("B718").Value = 3020587.8
$ws.Range("B719").Value = 3020587.8
